$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows of feed log data after the existing last row (161)
$ws.Range("A162").Value = 161
$ws.Range("B162").Value = 1
$ws.Range("C162").Value = "2024-06-18 09:13:46"
$ws.Range("D162").Value = 200
$ws.Range("E162").Value = 11

$ws.Range("A163").Value = 162
$ws.Range("B163").Value = 2
$ws.Range("C163").Value = "2024-06-18 09:13:46"
$ws.Range("D163").Value = 200
$ws.Range("E163").Value = 0
